# Fruta / hortaliza, semanal
# Insert a new weekly price-record row for "Uva" (Vega Monumental Concepción)
# above the existing row 42, pushing all subsequent records down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 42 (shifts old rows 42..102 down to 43..103).
$ws.Rows.Item(42).Insert()

# Populate the new row 42 with the new price observation.
$ws.Cells.Item(42, 1).Value = 11
$ws.Cells.Item(42, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value = "Bíobío"
$ws.Cells.Item(42, 4).Value = 44586
$ws.Cells.Item(42, 5).Value = 8
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100109
$ws.Cells.Item(42, 8).Value = "Uva"
$ws.Cells.Item(42, 9).Value = 100109001
$ws.Cells.Item(42, 10).Value = "Uva"
$ws.Cells.Item(42, 11).Value = "Superior Seedless"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 13).Value = 200
$ws.Cells.Item(42, 14).Value = 14000
$ws.Cells.Item(42, 15).Value = 15000
$ws.Cells.Item(42, 16).Value = 14500
$ws.Cells.Item(42, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(42, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(42, 19).Value = 967
$ws.Cells.Item(42, 20).Value = 15
